$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data block for the first date (rows 2-5) and the second date
# (rows 6-9) were swapped: the later week's figures (originally in rows 6-9)
# now belong to rows 2-5, and the earlier week's figures (originally in rows
# 2-5) now belong to rows 6-9. Columns affected: D (Fecha), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# R (Origen), S (Precio $/Kg).

$cols = @("D", "M", "N", "O", "P", "R", "S")
$rowPairs = @(@(2, 6), @(3, 7), @(4, 8), @(5, 9))

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
